$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 18:52"

# Apply updated COVID-19 country statistics (new data pull + re-sort by "Casos totales")
$ws.Range("B4").Value = 350013
$ws.Range("C4").Value = 13340
$ws.Range("D4").Value = 19247
$ws.Range("E4").Value = 320439
$ws.Range("B7").Value = 101178
$ws.Range("C7").Value = 1055
$ws.Range("E7").Value = 70866
$ws.Range("B12").Value = 30217
$ws.Range("C12").Value = 3148
$ws.Range("D12").Value = 1326
$ws.Range("E12").Value = 28242
$ws.Range("G12").Value = 75
$ws.Range("H12").Value = 649
$ws.Range("E13").Value = 13592
$ws.Range("G13").Value = 47
$ws.Range("H13").Value = 762
$ws.Range("A34").Value = "Ecuador"
$ws.Range("B34").Value = 3747
$ws.Range("C34").Value = 101
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = 3456
$ws.Range("F34").Value = 156
$ws.Range("G34").Value = 11
$ws.Range("H34").Value = 191
$ws.Range("A35").Value = "Pakistan"
$ws.Range("B35").Value = 3662
$ws.Range("C35").Value = 505
$ws.Range("D35").Value = 259
$ws.Range("E35").Value = 3351
$ws.Range("F35").Value = 17
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 52
$ws.Range("A36").Value = "Filipinas"
$ws.Range("B36").Value = 3660
$ws.Range("C36").Value = 414
$ws.Range("D36").Value = 73
$ws.Range("E36").Value = 3424
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 11
$ws.Range("H36").Value = 163
$ws.Range("A37").Value = "Japon"
$ws.Range("B37").Value = 3654
$ws.Range("D37").Value = 575
$ws.Range("E37").Value = 2994
$ws.Range("F37").Value = 69
$ws.Range("H37").Value = 85
$ws.Range("A76").Value = "Kazajistan"
$ws.Range("B76").Value = 651
$ws.Range("C76").Value = 67
$ws.Range("D76").Value = 45
$ws.Range("E76").Value = 600
$ws.Range("F76").Value = 16
$ws.Range("H76").Value = 6
$ws.Range("A77").Value = "Azerbaiyan"
$ws.Range("B77").Value = 641
$ws.Range("C77").Value = 57
$ws.Range("D77").Value = 44
$ws.Range("E77").Value = 590
$ws.Range("F77").Value = 11
$ws.Range("H77").Value = 7
$ws.Range("F102").Value = 2
$ws.Range("A156").Value = "Haiti"
$ws.Range("D156").Value = 0
$ws.Range("E156").Value = 23
$ws.Range("A157").Value = "Gabon"
$ws.Range("C157").Value = 3
$ws.Range("D157").Value = 1
$ws.Range("E157").Value = 22
$ws.Range("A158").Value = "Tanzania"
$ws.Range("B158").Value = 24
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 3
$ws.Range("E158").Value = 20
$ws.Range("G158").Value = 0
$ws.Range("A159").Value = "Benin"
$ws.Range("B159").Value = 23
$ws.Range("D159").Value = 5
$ws.Range("E159").Value = 17
$ws.Range("G159").Value = 1
$ws.Range("A160").Value = "Birmania"
$ws.Range("B160").Value = 22
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 21
$ws.Range("A164").Value = "Nueva Caledonia"
$ws.Range("D164").Value = 1
$ws.Range("H164").Value = 0
$ws.Range("A165").Value = "Libia"
$ws.Range("D165").Value = 0
$ws.Range("H165").Value = 1
$ws.Range("A171").Value = "Dominica"
$ws.Range("C171").Value = 0
$ws.Range("A172").Value = "Fiyi"
$ws.Range("C172").Value = 2
$ws.Range("D181").Value = 4
$ws.Range("E181").Value = 7
$ws.Range("A182").Value = "San Cristobal y Nieves"
$ws.Range("C182").Value = 0
$ws.Range("A183").Value = "Suazilandia"
$ws.Range("C183").Value = 1
$ws.Range("A184").Value = "Mozambique"
$ws.Range("D184").Value = 1
$ws.Range("H184").Value = 0
$ws.Range("A185").Value = "Surinam"
$ws.Range("D185").Value = 0
$ws.Range("H185").Value = 1
$ws.Range("A192").Value = "Somalia"
$ws.Range("A193").Value = "Belice"
$ws.Range("C193").Value = 2
$ws.Range("D193").Value = 0
$ws.Range("F193").Value = 1
$ws.Range("G193").Value = 1
$ws.Range("H193").Value = 1
$ws.Range("A194").Value = "San Vicente y las Granadinas"
$ws.Range("E194").Value = 6
$ws.Range("H194").Value = 0
$ws.Range("A195").Value = "Cabo Verde"
$ws.Range("B195").Value = 7
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 5
$ws.Range("H195").Value = 1
$ws.Range("A196").Value = "Sierra Leona"
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 6
$ws.Range("A199").Value = "San Bartolome"
$ws.Range("D199").Value = 1
$ws.Range("E199").Value = 5
$ws.Range("H199").Value = 0
$ws.Range("A200").Value = "Mauritania"
$ws.Range("B200").Value = 6
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 2
$ws.Range("E200").Value = 3
$ws.Range("F200").Value = 0
$ws.Range("H200").Value = 1
$ws.Range("A201").Value = "Malaui"
$ws.Range("C201").Value = 1
$ws.Range("E201").Value = 5
$ws.Range("F201").Value = 1
$ws.Range("H201").Value = 0
$ws.Range("A202").Value = "Islas Turcas y Caicos"
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("A208").Value = "Anguila"
$ws.Range("A209").Value = "Islas Virgenes Britanicas"
$ws.Range("A213").Value = "Sudan del Sur"
$ws.Range("A215").Value = "San Pedro y Miquelon"

Write-Host "Applied paises.xlsx update (countries & provincias Spain)"
